# Applies the "Vistas Orden" list restructuring on Hoja1 (sheet3):
#  - section header moves from C24 (paired with F24) to a single C23 cell
#  - the backlog list below it is reordered/extended from 21 to 24 items,
#    adding three brand new entries ("menu", "footer", "blank page") near
#    the top and two more ("caja", "cocina") further down, replacing the
#    old, more verbose entries that covered the same ground
#  - the old companion header cell F24 is removed entirely

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- 1) Prime the shared-string table so brand-new strings get created in
#        the same order the original author typed them in (this keeps the
#        saved xl/sharedStrings.xml ordering faithful to the source edit).
#        We use a scratch cell far away from any used range and clear it
#        again afterwards, so it leaves no visible trace.
$scratch = $ws.Range("Z1")
$scratch.Value = "blank page"
$scratch.Value = "menu"
$scratch.Value = "footer"
$scratch.Value = "caja"
$scratch.Value = "cocina"
$scratch.Clear()

# --- 2) Move/rewrite the section header.
# Copy the formatting of the old header cell (C24, style 19) onto the new
# header location (C23), then set its text.
$ws.Range("C24").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value = "Vistas Orden"

# Remove the old companion header that used to sit in column F.
$ws.Range("F24").Clear()

# --- 3) Rebuild the ordered backlog list in B24:C47.
# Give the whole block the regular data-row formatting (style 15) by
# copying from a cell that already carries it, then fill in the new
# sequence of items.
$ws.Range("C25").Copy()
$ws.Range("C24:C47").PasteSpecial(-4122)

$items = @(
    "menu",
    "footer",
    "blank page",
    "Login sistema",
    "menu principal",
    "Pedido selección mesa",
    "Pedido Selección Productos",
    "login usuario",
    "caja",
    "cocina",
    "pago de pedido",
    "detalle del pedido",
    "inventario",
    "pedidos de determinado usuario",
    "cuentas abiertas",
    "perfil de usuario",
    "mensajeria",
    "inicio administracion",
    "admin reportes",
    "admin crud usuarios",
    "admin crud mesas",
    "admin crud categorias",
    "admin crud productos",
    "admin atenciones"
)

$startRow = 24
for ($i = 0; $i -lt $items.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $i + 1
    $ws.Cells.Item($row, 3).Value = $items[$i]
}

# --- 4) Leave the selection the way the author left it.
$ws.Range("C26").Select()
